# Populate Sheet2 with a small data table (headers a/b/c plus numeric rows)
# and make Sheet2 the active sheet with D5 selected, matching the target
# workbook state.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "a"
$ws2.Range("B1").Value = "b"
$ws2.Range("C1").Value = "c"

# Data rows
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3

$ws2.Range("A3").Value = 4
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = 6

$ws2.Range("A4").Value = 7
$ws2.Range("B4").Value = 8
$ws2.Range("C4").Value = 9
$ws2.Range("D4").Value = 10

# Make Sheet2 the active sheet and select D5, as in the target workbook.
$ws2.Activate()
$ws2.Range("D5").Select()
